$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 303. This shifts the existing rows 303:386
# down to 305:388 (the previous rows 385:386 end up at 387:388), matching the
# "dimension" growth from A1:R386 to A1:R388. The original rows 301:302 stay
# put, but their old content needs to move down into the freshly inserted
# 303:304 so a brand-new observation can be written into 301:302.
$ws.Rows.Item(303).Insert()
$ws.Rows.Item(303).Insert()

$ws.Range("A301:R302").Copy()
$ws.Range("A303").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Populate rows 301 and 302 with the brand-new weekly
# observation (date 2021-09-24, serial 44463) for "Coliflor" Primera/Segunda.
$ws.Cells.Item(301, 1).Value = 8
$ws.Cells.Item(301, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(301, 3).Value = "Coquimbo"
$ws.Cells.Item(301, 4).Value = 44463
$ws.Cells.Item(301, 5).Value = 4
$ws.Cells.Item(301, 6).Value = 100112008
$ws.Cells.Item(301, 7).Value = "Coliflor"
$ws.Cells.Item(301, 8).Value = "Sin especificar"
$ws.Cells.Item(301, 9).Value = "Primera"
$ws.Cells.Item(301, 10).Value = 3600
$ws.Cells.Item(301, 11).Value = 600
$ws.Cells.Item(301, 12).Value = 700
$ws.Cells.Item(301, 13).Value = 650
$ws.Cells.Item(301, 14).Value = "`$/unidad"
$ws.Cells.Item(301, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(301, 16).Value = 650
$ws.Cells.Item(301, 17).Value = 1
$ws.Cells.Item(301, 18).Value = "Hortaliza"

$ws.Cells.Item(302, 1).Value = 8
$ws.Cells.Item(302, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(302, 3).Value = "Coquimbo"
$ws.Cells.Item(302, 4).Value = 44463
$ws.Cells.Item(302, 5).Value = 4
$ws.Cells.Item(302, 6).Value = 100112008
$ws.Cells.Item(302, 7).Value = "Coliflor"
$ws.Cells.Item(302, 8).Value = "Sin especificar"
$ws.Cells.Item(302, 9).Value = "Segunda"
$ws.Cells.Item(302, 10).Value = 1800
$ws.Cells.Item(302, 11).Value = 500
$ws.Cells.Item(302, 12).Value = 550
$ws.Cells.Item(302, 13).Value = 525
$ws.Cells.Item(302, 14).Value = "`$/unidad"
$ws.Cells.Item(302, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(302, 16).Value = 525
$ws.Cells.Item(302, 17).Value = 1
$ws.Cells.Item(302, 18).Value = "Hortaliza"
